$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.723.36'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.600.76'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = "'211.61"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = "'0.513"
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D8").Value = "'0.0618"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").Value = "'19.67"
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").Value = '1.827.08'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '1.596.59'
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").Value = "'0.523"
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").Value = "'65.05"
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = "'210.13"
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").Value = "'7.13"
$ws.Range("E20").Value = '  +1.54%  '
$ws.Range("D21").Value = "'4.28"
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("E22").Value = '  -3.09%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = "'143.55"
$ws.Range("E24").Value = '  -0.53%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").Value = "'15.34"
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").Value = '1.287.60'
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("E36").Value = '  -2.84%  '
$ws.Range("E37").Value = '  +10.52%  '
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("E40").Value = '  -2.04%  '
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").Value = "'0.783"
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").Value = "'62.88"
$ws.Range("E43").Value = '  -1.15%  '
$ws.Range("D44").Value = '1.737.54'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = "'90.46"
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("E48").Value = '  +1.55%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("D50").Value = "'7.46"
$ws.Range("E50").Value = '  +0.55%  '
